$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3 value from "ahti" to "samuli"
$ws.Range("A3").Value = "samuli"

# Add new row 4 with "markus"
$ws.Range("A4").Value = "markus"

# Add numeric values in column B
$ws.Range("B1").Value = 2
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 7
$ws.Range("B4").Value = 10

# Update selection to match target (F8)
$ws.Range("F8").Select()
